# Applies the "Updated symbol list" data refresh to the crypto price sheet.
# Columns: D = Price, E = Volume(1h), G = Hora (hour)
# All values are written as literal text to match the original inlineStr
# cell contents exactly (preserving trailing zeros, percent signs, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    # Force the cell to stay text (avoid Excel auto-converting numeric-looking
    # strings like "5.620" or "1.60%" into numbers and losing precision),
    # then restore the default "Normal" style so no stray formatting is left
    # behind (NumberFormat = "@" alone would leave a text-format style).
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$updates = @(
    @{Row=2; D="329.11"; E="1.60%"; G="6"},
    @{Row=3; D="41.25"; E="3.84%"; G="6"},
    @{Row=4; D="5.620"; E="-4.44%"; G="6"},
    @{Row=5; D="0.08171"; E="1.79%"; G="6"},
    @{Row=6; D="2.023"; E="5.12%"; G="6"},
    @{Row=7; D="8.742"; E="0.81%"; G="6"},
    @{Row=8; D="4.533"; E="-0.88%"; G="6"},
    @{Row=9; D="2.987"; E="1.35%"; G="6"},
    @{Row=10; D="0.9176"; E="-1.51%"; G="6"},
    @{Row=11; D="0.1264"; E="0.45%"; G="6"},
    @{Row=12; D="0.1950"; E="-1.01%"; G="6"},
    @{Row=13; D="0.09316"; E="1.11%"; G="6"},
    @{Row=14; D="0.03742"; E="5.32%"; G="6"},
    @{Row=15; D="0.1059"; E="1.21%"; G="6"},
    @{Row=16; D="0.001307"; E="1.24%"; G="6"},
    @{Row=17; D="0.006158"; E="-0.38%"; G="6"},
    @{Row=18; D="3.437"; E="2.66%"; G="6"},
    @{Row=19; E="-1.47%"; G="6"},
    @{Row=20; D="8.306"; E="-5.02%"; G="6"},
    @{Row=21; D="0.1380"; E="-2.74%"; G="6"},
    @{Row=22; D="0.2392"; E="-2.27%"; G="6"},
    @{Row=23; D="0.04426"; E="0.12%"; G="6"},
    @{Row=24; E="0.01%"; G="6"},
    @{Row=25; D="0.004283"; E="-2.59%"; G="6"},
    @{Row=26; D="0.0001183"; E="3.74%"; G="6"},
    @{Row=27; G="6"},
    @{Row=28; G="6"},
    @{Row=29; G="6"},
    @{Row=30; G="6"},
    @{Row=31; G="6"},
    @{Row=32; G="6"},
    @{Row=33; G="6"},
    @{Row=34; G="6"},
    @{Row=35; G="6"},
    @{Row=36; G="6"},
    @{Row=37; G="6"},
    @{Row=38; G="6"},
    @{Row=39; D="0.02743"; E="9.72%"; G="6"},
    @{Row=40; D="0.05401"; E="2.12%"; G="6"},
    @{Row=41; D="0.007660"; E="3.42%"; G="6"},
    @{Row=42; D="0.1413"; E="0.55%"; G="6"},
    @{Row=43; D="0.008999"; E="-6.27%"; G="6"},
    @{Row=44; D="0.002135"; E="0.86%"; G="6"},
    @{Row=45; D="0.01129"; G="6"},
    @{Row=46; D="0.00006876"; E="2.06%"; G="6"},
    @{Row=47; D="0.00000000752"; E="0.22%"; G="6"},
    @{Row=48; D="0.003583"; E="19.41%"; G="6"},
    @{Row=49; D="0.002284"; E="60.54%"; G="6"},
    @{Row=50; D="0.00002105"; E="0.22%"; G="6"},
    @{Row=51; D="0.0002004"; E="0.22%"; G="6"}

)

foreach ($item in $updates) {
    $row = $item.Row
    if ($item.ContainsKey("D")) {
        Set-TextCell $ws.Cells.Item($row, 4) $item.D
    }
    if ($item.ContainsKey("E")) {
        Set-TextCell $ws.Cells.Item($row, 5) $item.E
    }
    if ($item.ContainsKey("G")) {
        Set-TextCell $ws.Cells.Item($row, 7) $item.G
    }
}
